# Generate Report for Handoff
# Update "Latest Handoff Date/Datetime" for the a01f0b9d-d8c8-4a3d-aa82-8defee31c6ee
# file (row 5 of each data sheet) to reflect a fresh handoff report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: column D = "Latest Handoff Date"
$wsOverview.Range("D5").Value = "2016-03-23 09:17:33"

# zh-cn sheet: column E = "Latest Handoff Datetime"
$wsZhCn.Range("E5").Value = "2016-03-23 09:17:24"

# de-de sheet: column E = "Latest Handoff Datetime"
$wsDeDe.Range("E5").Value = "2016-03-23 09:17:33"
